$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'The Lord of the Rings: The War of the Rohirrim'
$ws.Range("C2").Value = "'6.6"

$ws.Range("A3").Value = 'Gladiator II'
$ws.Range("C3").Value = "'6.8"

$ws.Range("A4").Value = 'The Order'
$ws.Range("B4").Value = 'Crime'
$ws.Range("C4").Value = "'7"

$ws.Range("A5").Value = 'The Six Triple Eight'
$ws.Range("B5").Value = 'Drama'
$ws.Range("C5").Value = "'6.5"

$ws.Range("A6").Value = 'Bird'
$ws.Range("B6").Value = 'Drama'
$ws.Range("C6").Value = "'7.2"

$ws.Range("A7").Value = 'Small Things Like These'
$ws.Range("B7").Value = 'Drama'

$ws.Range("A8").Value = 'Anora'
$ws.Range("B8").Value = 'Comedy'
$ws.Range("C8").Value = "'8"

$ws.Range("A9").Value = 'Dirty Angels'
$ws.Range("B9").Value = 'Action'
$ws.Range("C9").Value = "'4.1"

$ws.Range("A10").Value = 'Omni Loop'
$ws.Range("B10").Value = 'Drama'
$ws.Range("C10").Value = "'5.5"

$ws.Range("A11").Value = 'DragonHeart'
$ws.Range("C11").Value = "'6.4"
$ws.Range("D11").Value = "'1996"

$ws.Range("A12").Value = 'Megalopolis'
$ws.Range("B12").Value = 'Drama'
$ws.Range("C12").Value = "'4.8"

$ws.Range("A13").Value = 'Hush'
$ws.Range("B13").Value = 'Action'
$ws.Range("C13").Value = "'6.6"
$ws.Range("D13").Value = "'2016"

$ws.Range("A14").Value = 'Candyman'
$ws.Range("C14").Value = "'6.7"
$ws.Range("D14").Value = "'1992"

$ws.Range("A15").Value = 'Teenage Mutant Ninja Turtles: Out of the Shadows'
$ws.Range("C15").Value = "'5.9"
$ws.Range("D15").Value = "'2016"

$ws.Range("A16").Value = 'Demolition Man'
$ws.Range("C16").Value = "'6.7"
$ws.Range("D16").Value = "'1993"

$ws.Range("A17").Value = '[ES] The Settlers'
$ws.Range("B17").Value = 'Action'
$ws.Range("C17").Value = "'6.9"
$ws.Range("D17").Value = "'2023"

$ws.Range("A18").Value = 'Twin Peaks: Fire Walk with Me'
$ws.Range("C18").Value = "'7.3"
$ws.Range("D18").Value = "'1992"

$ws.Range("A19").Value = 'The Last House on the Left'
$ws.Range("C19").Value = "'6.5"
$ws.Range("D19").Value = "'2009"

$ws.Range("A20").Value = 'Ted'
$ws.Range("B20").Value = 'Action'
$ws.Range("C20").Value = "'6.9"
$ws.Range("D20").Value = "'2012"

$ws.Range("A21").Value = 'Starve Acre'
$ws.Range("C21").Value = "'5.4"
